$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 102:103, shifting the existing rows (102-170) down to (104-172).
$ws.Rows("102:103").Insert(-4121)

# The block that used to be rows 102:103 is now at rows 104:105. Duplicate that
# data back into the newly blank rows 102:103 so every column (A-R) is populated
# exactly like the rest of the sheet, then we will overwrite the handful of
# cells (D, K, L, M, P) that differ for the new rows.
$src = $ws.Range("A104:R105")
$dst = $ws.Range("A102:R103")
$src.Copy($dst)

# Row 102: new date + updated min/max/avg/price values.
$ws.Range("D102").Value = 44596
$ws.Range("K102").Value = 450
$ws.Range("L102").Value = 450
$ws.Range("M102").Value = 450
$ws.Range("P102").Value = 450

# Row 103: new date + updated min/max/avg/price values.
$ws.Range("D103").Value = 44596
$ws.Range("K103").Value = 350
$ws.Range("L103").Value = 350
$ws.Range("M103").Value = 350
$ws.Range("P103").Value = 350
